$d = $word.ActiveDocument

# Trim the trailing sentence from the "message file not in a 26-bit chunks
# format" bullet so it ends right after "It won't be able to be read. "
# Re-supplying the whole sentence (rather than just deleting the tail) lets
# Word re-flow the paragraph's runs into a single run, matching how the
# author's edit collapsed the previously-split runs.
$oldText = "In case the message file is not in a 26-bit chunks format. It won" + [char]0x2019 + "t be able to be read. Same thing would happen in case the message would include chars different from " + [char]0x2018 + "1" + [char]0x2019 + " or " + [char]0x2018 + "0" + [char]0x2019 + "."
$newText = "In case the message file is not in a 26-bit chunks format. It won" + [char]0x2019 + "t be able to be read. "

$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null
